$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.172.74"
$ws.Range("E2").Value = "  +5.34%  "
$ws.Range("D3").Value = "2.999.40"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'580.94"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "'162.64"
$ws.Range("E6").Value = "  +12.14%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.99%  "
$ws.Range("D9").Value = "2.995.96"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("E10").Value = "  -5.16%  "
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("E12").Value = "  +4.55%  "
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  +5.31%  "
$ws.Range("D14").Value = "'34.55"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "66.165.79"
$ws.Range("E16").Value = "  +5.39%  "
$ws.Range("D17").Value = "3.497.63"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "'6.91"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").Value = "3.001.29"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").Value = "'453.78"
$ws.Range("E20").Value = "  +5.13%  "
$ws.Range("D21").Value = "'13.83"
$ws.Range("E21").Value = "  +5.19%  "
$ws.Range("D22").Value = "'0.685"
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("D23").Value = "'7.33"
$ws.Range("E23").Value = "  +5.81%  "
$ws.Range("D24").Value = "'82.30"
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  +13.50%  "
$ws.Range("D26").Value = "'12.25"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").Value = "'10.37"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'8.13"
$ws.Range("E29").Value = "  +12.91%  "
$ws.Range("E30").Value = "  +17.66%  "
$ws.Range("D31").Value = "'2.62"
$ws.Range("E31").Value = "  +4.80%  "
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").Value = "'27.30"
$ws.Range("E33").Value = "  +5.23%  "
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").Value = "'5.81"
$ws.Range("E37").Value = "  +7.22%  "
$ws.Range("D38").Value = "'2.06"
$ws.Range("E38").Value = "  +7.21%  "
$ws.Range("D39").Value = "'49.66"
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").Value = "'2.96"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "'0.309"
$ws.Range("E41").Value = "  +15.30%  "
$ws.Range("E42").Value = "  +6.66%  "
$ws.Range("D43").Value = "'43.87"
$ws.Range("E43").Value = "  +7.49%  "
$ws.Range("D44").Value = "'8.42"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("D45").Value = "'399.74"
$ws.Range("E45").Value = "  +12.17%  "
$ws.Range("D46").Value = "'0.0357"
$ws.Range("E46").Value = "  +4.84%  "
$ws.Range("D47").Value = "2.787.79"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("D48").Value = "'133.70"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'23.82"
$ws.Range("E50").Value = "  +10.31%  "
$ws.Range("E51").Value = "  +3.26%  "
